$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.065160274505615
$ws.Range("B1").Value = 2.630269050598145
$ws.Range("C1").Value = 2.015894174575806
$ws.Range("D1").Value = 1.798981785774231
$ws.Range("E1").Value = 1.611698269844055
